# Added Taser and fixed it.
# Adds two new rows of translation data (IDs, English, Ukrainian) for the
# "Ammo" and "Infinite Taser" pickup strings, appending four new shared
# strings and leaving column C (Spanish) empty for both rows, matching the
# upstream commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A125").Value = "Ammo"
$ws.Range("B125").Value = "Ammo"
$ws.Range("D125").Value = "Патрони"

$ws.Range("A126").Value = "InfiniteTaser"
$ws.Range("B126").Value = "Infinite Taser"
$ws.Range("D126").Value = "Бескінечний тазер"

# Mirror the author's final cursor position/selection from the diff.
$ws.Range("D127").Select()
